# Trading update: 2026-02-17 15:30:25
# Appends the newly-opened trade (row 46) to both the "All Trades" and the
# "MarketMaking" worksheets, which are kept in sync with identical data.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Trade #
    $ws.Range("A46").Value = 45

    # Date - force text so it is not auto-converted into a date serial value
    $ws.Range("B46").NumberFormat = "@"
    $ws.Range("B46").Value = "2026-02-17"
    $ws.Range("B46").Style = "Normal"

    # Time
    $ws.Range("C46").Value = "15:30:24"

    # Strategy / Side
    $ws.Range("D46").Value = "MarketMaking"
    $ws.Range("E46").Value = "DOWN"

    # Entry Price
    $ws.Range("F46").Value = 0.378323

    # Exit Price - trade is still open, so this stays blank
    $ws.Range("G46").NumberFormat = "@"
    $ws.Range("G46").Value = ""
    $ws.Range("G46").Style = "Normal"

    # Status
    $ws.Range("H46").Value = "OPEN"

    # P&L %, P&L $
    $ws.Range("I46").Value = 0
    $ws.Range("J46").Value = 0

    # Capital After
    $ws.Range("K46").Value = 100.7498800637372

    # Entry/Exit Slippage (bps)
    $ws.Range("L46").Value = 0
    $ws.Range("M46").Value = 0

    # Confidence
    $ws.Range("N46").Value = 0.6

    # Entry Reason
    $ws.Range("O46").Value = "Normal spread capture: 19600 bps"

    # Exit Reason - blank, trade still open
    $ws.Range("P46").NumberFormat = "@"
    $ws.Range("P46").Value = ""
    $ws.Range("P46").Style = "Normal"

    # Duration (min)
    $ws.Range("Q46").Value = 0
}
